$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 20).Value = 0
$ws.Cells.Item(2, 21).Value = 0
$ws.Cells.Item(2, 22).Value = 0.0066
$ws.Cells.Item(2, 23).Value = 0
$ws.Cells.Item(2, 24).Value = 0
$ws.Cells.Item(3, 20).Value = 0
$ws.Cells.Item(3, 21).Value = 0
$ws.Cells.Item(3, 22).Value = 0.0066
$ws.Cells.Item(3, 24).Value = 0
$ws.Cells.Item(4, 20).Value = 0
$ws.Cells.Item(4, 21).Value = 0
$ws.Cells.Item(4, 22).Value = 0
$ws.Cells.Item(4, 23).Value = 0
$ws.Cells.Item(4, 24).Value = 0
$ws.Cells.Item(5, 20).Value = 0
$ws.Cells.Item(5, 21).Value = 0
$ws.Cells.Item(5, 22).Value = 0
$ws.Cells.Item(5, 23).Value = 0
$ws.Cells.Item(6, 20).Value = 0
$ws.Cells.Item(6, 21).Value = 0.0066
$ws.Cells.Item(6, 22).Value = 0
$ws.Cells.Item(6, 23).Value = 0
$ws.Cells.Item(6, 24).Value = 0
$ws.Cells.Item(7, 20).Value = 0
$ws.Cells.Item(7, 21).Value = 0
$ws.Cells.Item(7, 22).Value = 0
$ws.Cells.Item(7, 23).Value = 0
$ws.Cells.Item(7, 24).Value = 0.0066
$ws.Cells.Item(8, 20).Value = 0
$ws.Cells.Item(8, 21).Value = 0
$ws.Cells.Item(8, 22).Value = 0.0066
$ws.Cells.Item(8, 23).Value = 0
$ws.Cells.Item(8, 24).Value = 0
$ws.Cells.Item(9, 20).Value = 0
$ws.Cells.Item(9, 21).Value = 0
$ws.Cells.Item(9, 23).Value = 0
$ws.Cells.Item(9, 24).Value = 0
$ws.Cells.Item(10, 20).Value = 0
$ws.Cells.Item(10, 21).Value = 0
$ws.Cells.Item(10, 22).Value = 0
$ws.Cells.Item(10, 23).Value = 0
$ws.Cells.Item(10, 24).Value = 0
$ws.Cells.Item(11, 20).Value = 0
$ws.Cells.Item(11, 21).Value = 0
$ws.Cells.Item(11, 22).Value = 0.0066
$ws.Cells.Item(11, 23).Value = 0
$ws.Cells.Item(11, 24).Value = 0.0066
$ws.Cells.Item(12, 20).Value = 0
$ws.Cells.Item(12, 21).Value = 0
$ws.Cells.Item(12, 22).Value = 0.0264
$ws.Cells.Item(12, 23).Value = 0.0066
$ws.Cells.Item(12, 24).Value = 0
$ws.Cells.Item(13, 20).Value = 0.0132
$ws.Cells.Item(13, 21).Value = 0.0132
$ws.Cells.Item(13, 22).Value = 0.0198
$ws.Cells.Item(13, 23).Value = 0
$ws.Cells.Item(13, 24).Value = 0
$ws.Cells.Item(14, 20).Value = 0.0066
$ws.Cells.Item(14, 21).Value = 0.0066
$ws.Cells.Item(14, 22).Value = 0.0066
$ws.Cells.Item(14, 23).Value = 0.0066
$ws.Cells.Item(14, 24).Value = 0
$ws.Cells.Item(15, 20).Value = 0
$ws.Cells.Item(15, 21).Value = 0.0198
$ws.Cells.Item(15, 22).Value = 0.0066
$ws.Cells.Item(15, 23).Value = 0.0066
$ws.Cells.Item(15, 24).Value = 0.0132
$ws.Cells.Item(16, 20).Value = 0
$ws.Cells.Item(16, 21).Value = 0.0066
$ws.Cells.Item(16, 22).Value = 0.0066
$ws.Cells.Item(16, 23).Value = 0.0066
$ws.Cells.Item(16, 24).Value = 0
$ws.Cells.Item(17, 20).Value = 0
$ws.Cells.Item(17, 21).Value = 0.0066
$ws.Cells.Item(17, 22).Value = 0.0198
$ws.Cells.Item(17, 23).Value = 0
$ws.Cells.Item(17, 24).Value = 0.0132
$ws.Cells.Item(18, 20).Value = 0
$ws.Cells.Item(18, 21).Value = 0
$ws.Cells.Item(18, 22).Value = 0.0198
$ws.Cells.Item(18, 23).Value = 0
$ws.Cells.Item(18, 24).Value = 0.0198
$ws.Cells.Item(19, 20).Value = 0.0066
$ws.Cells.Item(19, 21).Value = 0.0066
$ws.Cells.Item(19, 22).Value = 0.0066
$ws.Cells.Item(19, 23).Value = 0
$ws.Cells.Item(19, 24).Value = 0.0264
$ws.Cells.Item(20, 20).Value = 0.0132
$ws.Cells.Item(20, 21).Value = 0.0132
$ws.Cells.Item(20, 22).Value = 0.0198
$ws.Cells.Item(20, 23).Value = 0.0066
$ws.Cells.Item(20, 24).Value = 0.0198
$ws.Cells.Item(21, 20).Value = 0.0066
$ws.Cells.Item(21, 21).Value = 0.0066
$ws.Cells.Item(21, 22).Value = 0.0066
$ws.Cells.Item(21, 23).Value = 0.0066
$ws.Cells.Item(21, 24).Value = 0.0132
$ws.Cells.Item(22, 20).Value = 0.0066
$ws.Cells.Item(22, 21).Value = 0.0066
$ws.Cells.Item(22, 22).Value = 0.0132
$ws.Cells.Item(22, 23).Value = 0.0066
$ws.Cells.Item(22, 24).Value = 0.0066
$ws.Cells.Item(23, 20).Value = 0.0066
$ws.Cells.Item(23, 21).Value = 0.0132
$ws.Cells.Item(23, 22).Value = 0.0264
$ws.Cells.Item(23, 23).Value = 0
$ws.Cells.Item(23, 24).Value = 0.0066
$ws.Cells.Item(24, 20).Value = 0
$ws.Cells.Item(24, 21).Value = 0.0066
$ws.Cells.Item(24, 22).Value = 0.0264
$ws.Cells.Item(24, 23).Value = 0.0066
$ws.Cells.Item(24, 24).Value = 0
$ws.Cells.Item(25, 20).Value = 0
$ws.Cells.Item(25, 21).Value = 0
$ws.Cells.Item(25, 22).Value = 0.0066
$ws.Cells.Item(25, 23).Value = 0.0066
$ws.Cells.Item(25, 24).Value = 0
